$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the D-column cells whose new values are ambiguous (look like numbers)
# as Text so Excel stores them as literal strings instead of coercing to numbers,
# matching the source data (t="inlineStr") without touching any other cell style.
$textCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D17', 'D18', 'D19', 'D21', 'D22', 'D23', 'D25', 'D26', 'D27', 'D28', 'D29', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '26.590.26'
$ws.Range('E2').Value = '  +6.98%  '

# Row 3
$ws.Range('D3').Value = '1.721.88'
$ws.Range('E3').Value = '  +3.39%  '

# Row 4
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.30%  '

# Row 5
$ws.Range('D5').Value = '333.82'
$ws.Range('E5').Value = '  +1.53%  '

# Row 6
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').Value = '  -0.13%  '

# Row 7
$ws.Range('D7').Value = '0.3709'
$ws.Range('E7').Value = '  +1.75%  '

# Row 8
$ws.Range('D8').Value = '48.25'
$ws.Range('E8').Value = '  +1.79%  '

# Row 9
$ws.Range('D9').Value = '0.3359'
$ws.Range('E9').Value = '  +2.81%  '

# Row 10
$ws.Range('D10').Value = '1.186'
$ws.Range('E10').Value = '  +4.20%  '

# Row 11
$ws.Range('D11').Value = '0.07400'
$ws.Range('E11').Value = '  +4.42%  '

# Row 12
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.08%  '

# Row 13
$ws.Range('D13').Value = '6.379'
$ws.Range('E13').Value = '  +4.76%  '

# Row 14
$ws.Range('D14').Value = '20.10'
$ws.Range('E14').Value = '  +2.98%  '

# Row 15
$ws.Range('D15').Value = '7.041'
$ws.Range('E15').Value = '  +6.49%  '

# Row 16
$ws.Range('D16').Value = '1.713.56'
$ws.Range('E16').Value = '  +2.80%  '

# Row 17
$ws.Range('D17').Value = '0.00001070'
$ws.Range('E17').Value = '  +1.97%  '

# Row 18
$ws.Range('D18').Value = '0.06631'
$ws.Range('E18').Value = '  -0.11%  '

# Row 19
$ws.Range('D19').Value = '81.98'
$ws.Range('E19').Value = '  +4.40%  '

# Row 20
$ws.Range('E20').Value = '  +0.04%  '

# Row 21
$ws.Range('D21').Value = '16.54'
$ws.Range('E21').Value = '  +4.60%  '

# Row 22
$ws.Range('D22').Value = '6.142'
$ws.Range('E22').Value = '  +3.63%  '

# Row 23
$ws.Range('D23').Value = '12.77'
$ws.Range('E23').Value = '  +1.66%  '

# Row 24
$ws.Range('D24').Value = '26.514.89'
$ws.Range('E24').Value = '  +6.69%  '

# Row 25
$ws.Range('D25').Value = '2.428'
$ws.Range('E25').Value = '  -1.56%  '

# Row 26
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').Value = '1.400'
$ws.Range('E26').Value = '  +19.60%  '

# Row 27
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.399'
$ws.Range('E27').Value = '  -1.46%  '

# Row 28
$ws.Range('D28').Value = '151.82'
$ws.Range('E28').Value = '  +1.66%  '

# Row 29
$ws.Range('D29').Value = '19.38'
$ws.Range('E29').Value = '  +3.64%  '

# Row 30
$ws.Range('D30').Value = '1.911.69'
$ws.Range('E30').Value = '  +3.36%  '

# Row 31
$ws.Range('D31').Value = '131.08'
$ws.Range('E31').Value = '  +3.86%  '

# Row 32
$ws.Range('D32').Value = '4.115'
$ws.Range('E32').Value = '  +1.10%  '

# Row 33
$ws.Range('D33').Value = '5.945'
$ws.Range('E33').Value = '  +4.36%  '

# Row 34
$ws.Range('D34').Value = '0.08608'
$ws.Range('E34').Value = '  +1.21%  '

# Row 35
$ws.Range('D35').Value = '1.696'
$ws.Range('E35').Value = '  +2.58%  '

# Row 36
$ws.Range('D36').Value = '12.75'
$ws.Range('E36').Value = '  +4.56%  '

# Row 37
$ws.Range('D37').Value = '5.362'
$ws.Range('E37').Value = '  +3.87%  '

# Row 38
$ws.Range('D38').Value = '0.02326'
$ws.Range('E38').Value = '  +2.13%  '

# Row 39
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.06213'
$ws.Range('E39').Value = '  -0.21%  '

# Row 40
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').Value = '0.2153'
$ws.Range('E40').Value = '  +3.34%  '

# Row 41
$ws.Range('D41').Value = '8.445'
$ws.Range('E41').Value = '  +2.54%  '

# Row 42
$ws.Range('D42').Value = '1.221'
$ws.Range('E42').Value = '  -2.23%  '

# Row 43
$ws.Range('D43').Value = '0.6187'
$ws.Range('E43').Value = '  +4.11%  '

# Row 44
$ws.Range('D44').Value = '14.16'
$ws.Range('E44').Value = '  +5.17%  '

# Row 45
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.02%  '

# Row 46
$ws.Range('D46').Value = '3.901'
$ws.Range('E46').Value = '  +1.44%  '

# Row 47
$ws.Range('D47').Value = '0.5984'
$ws.Range('E47').Value = '  +5.78%  '

# Row 48
$ws.Range('D48').Value = '129.24'
$ws.Range('E48').Value = '  +3.00%  '

# Row 49
$ws.Range('D49').Value = '2.040'
$ws.Range('E49').Value = '  +4.31%  '

# Row 50
$ws.Range('D50').Value = '0.07166'
$ws.Range('E50').Value = '  +2.54%  '

# Row 51
$ws.Range('D51').Value = '76.84'
$ws.Range('E51').Value = '  +2.03%  '
